$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells D1:G1 to be prefixed with "label:"
$ws.Range("D1").Value = "label:role"
$ws.Range("E1").Value = "label:app"
$ws.Range("F1").Value = "label:env"
$ws.Range("G1").Value = "label:loc"

[void]$ws.Range("G8").Select()
